# Updated legacy GSC export data.
#
# The "Chart" sheet had two obsolete leading date rows (2025-09-28 and the
# blank-label row that immediately followed it, then 2025-09-29) that need
# to be dropped so the export starts at 2025-09-30. Deleting the two rows
# shifts every following row up by two and shrinks the used range from
# A1:D88 down to A1:D86 (matching the refreshed export).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 is always the next "oldest" row once the previous one is removed,
# so deleting it twice drops the 2025-09-28 and 2025-09-29 rows in turn.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
